$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "36.577.45"
$ws.Range("E2").Value = "  +0.42%  "
$ws.Range("D3").Value = "1.959.76"
$ws.Range("E3").Value = "  +0.94%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "243.53"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.14%  "
$ws.Range("E6").Value = "  +0.59%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "60.53"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +5.77%  "
$ws.Range("E8").Value = "  -0.03%  "
$ws.Range("E9").Value = "  +3.94%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0789"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -6.35%  "
$ws.Range("E11").Value = "  +0.16%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "14.21"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +5.53%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "21.78"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +2.45%  "
$ws.Range("D14").Value = "2.246.34"
$ws.Range("E14").Value = "  +0.79%  "
$ws.Range("E15").Value = "  +2.27%  "
$ws.Range("E16").Value = "  +2.03%  "
$ws.Range("D17").Value = "1.963.23"
$ws.Range("E17").Value = "  +0.53%  "
$ws.Range("D18").Value = "36.498.35"
$ws.Range("E18").Value = "  +0.27%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "69.61"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.45%  "
$ws.Range("D20").Value = "0.0₃0852"
$ws.Range("E20").Value = "  -1.23%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "229.50"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.42%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.07"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.44%  "
$ws.Range("E23").Value = "  +0.07%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.44"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +3.52%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.35"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +2.36%  "
$ws.Range("E26").Value = "  +5.32%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.20"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.12%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "161.12"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.43%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "19.32"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.74%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.33"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +21.53%  "
$ws.Range("E31").Value = "  +1.12%  "
$ws.Range("E32").Value = "  +4.12%  "
$ws.Range("E33").Value = "  -0.38%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.43"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +5.95%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.46"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +9.11%  "
$ws.Range("E36").Value = "  +0.01%  "
$ws.Range("E37").Value = "  +4.56%  "
$ws.Range("E38").Value = "  -1.08%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.48"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -12.42%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0968"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -2.09%  "
$ws.Range("E41").Value = "  +0.58%  "
$ws.Range("E42").Value = "  +1.14%  "
$ws.Range("E43").Value = "  +0.04%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "15.85"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.50%  "
$ws.Range("D45").Value = "1.363.75"
$ws.Range("E45").Value = "  +1.71%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "88.72"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +2.42%  "
$ws.Range("E47").Value = "  +0.02%  "
$ws.Range("E48").Value = "  -1.20%  "
$ws.Range("E49").Value = "  +0.12%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "45.79"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +5.90%  "
$ws.Range("D51").Value = "2.137.89"
$ws.Range("E51").Value = "  +0.88%  "
